$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '20.543.66'
$ws.Range('E2').Value = '  +1.83%  '
$ws.Range('D3').Value = '1.471.50'
$ws.Range('E3').Value = '  +2.79%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.008'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.81%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.9577'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -3.80%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '276.22'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.52%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3646'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -1.64%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3059'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -3.32%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '39.81'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -1.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.058'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.17%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06629'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.65%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.001'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.45%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.468'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.73%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.09'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.82%  '
$ws.Range('E15').Value = '  -1.02%  '
$ws.Range('E16').Value = '  -0.47%  '
$ws.Range('D17').Value = '1.472.39'
$ws.Range('E17').Value = '  +3.08%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.05916'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +2.87%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9633'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -3.18%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.12'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -4.26%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.453'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -3.44%  '
$ws.Range('E22').Value = '  -2.70%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.07'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.63%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.244'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.51%  '
$ws.Range('D25').Value = '20.588.07'
$ws.Range('E25').Value = '  +1.92%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '140.33'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +3.75%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.122'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -8.24%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.22'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.42%  '
$ws.Range('D29').Value = '1.629.80'
$ws.Range('E29').Value = '  +2.66%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '113.78'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.86%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.995'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.03%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.969'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -6.83%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.8113'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -4.16%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.07951'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +1.42%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.539'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +3.24%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.212'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +8.91%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.05823'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.18%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.703'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -4.94%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02036'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.79%  '
$ws.Range('E40').Value = '  -3.91%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9588'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -3.64%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '7.590'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -4.91%  '
$ws.Range('E43').Value = '  -0.32%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5276'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.507'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.54%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.06'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -2.96%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '117.81'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5177'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -2.24%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.784'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.92%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06461'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +3.17%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9911'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.55%  '
